$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 2).Value = 5143697
$ws.Cells.Item(3, 7).Value = 'FC Nordsjaelland'
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 11).Value = 3.3
$ws.Cells.Item(3, 12).Value = 3.6
$ws.Cells.Item(3, 13).Value = 2
$ws.Cells.Item(3, 14).Value = 4.2
$ws.Cells.Item(3, 15).Value = 3.75
$ws.Cells.Item(3, 16).Value = 1.8
$ws.Cells.Item(3, 17).Value = 0.75
$ws.Cells.Item(3, 18).Value = 1.8
$ws.Cells.Item(3, 19).Value = 2.05
$ws.Cells.Item(3, 21).Value = 1.975
$ws.Cells.Item(3, 22).Value = 1.875
$ws.Cells.Item(3, 24).Value = 2.75
$ws.Cells.Item(3, 26).Value = 0.8
$ws.Cells.Item(3, 27).Value = -1
$ws.Cells.Item(3, 29).Value = 0.875
$ws.Cells.Item(4, 2).Value = 5143696
$ws.Cells.Item(4, 7).Value = 'Randers FC'
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 2.2
$ws.Cells.Item(4, 12).Value = 3.4
$ws.Cells.Item(4, 13).Value = 3
$ws.Cells.Item(4, 14).Value = 2.1
$ws.Cells.Item(4, 15).Value = 3.5
$ws.Cells.Item(4, 16).Value = 3.3
$ws.Cells.Item(4, 17).Value = -0.25
$ws.Cells.Item(4, 18).Value = 1.85
$ws.Cells.Item(4, 19).Value = 2
$ws.Cells.Item(4, 21).Value = 1.85
$ws.Cells.Item(4, 22).Value = 2
$ws.Cells.Item(4, 24).Value = 2.5
$ws.Cells.Item(4, 26).Value = -0.5
$ws.Cells.Item(4, 27).Value = 0.5
$ws.Cells.Item(4, 29).Value = 1
$ws.Cells.Item(8, 7).Value = 'Odense BK'
$ws.Cells.Item(9, 7).Value = 'Lyngby'
$ws.Cells.Item(17, 6).Value = 'Lyngby'
$ws.Cells.Item(18, 7).Value = 'Odense BK'
$ws.Cells.Item(20, 6).Value = 'Odense BK'
$ws.Cells.Item(22, 7).Value = 'Lyngby'
$ws.Cells.Item(26, 2).Value = 5143711
$ws.Cells.Item(26, 6).Value = 'Silkeborg IF'
$ws.Cells.Item(26, 7).Value = 'Midtjylland'
$ws.Cells.Item(26, 8).Value = 3
$ws.Cells.Item(26, 9).Value = 3
$ws.Cells.Item(26, 11).Value = 2.25
$ws.Cells.Item(26, 12).Value = 3.6
$ws.Cells.Item(26, 13).Value = 2.9
$ws.Cells.Item(26, 14).Value = 2.25
$ws.Cells.Item(26, 15).Value = 3.75
$ws.Cells.Item(26, 16).Value = 2.875
$ws.Cells.Item(26, 18).Value = 1.975
$ws.Cells.Item(26, 19).Value = 1.875
$ws.Cells.Item(26, 20).Value = 2.75
$ws.Cells.Item(26, 21).Value = 1.95
$ws.Cells.Item(26, 22).Value = 1.9
$ws.Cells.Item(26, 24).Value = 2.75
$ws.Cells.Item(26, 27).Value = 0.4375
$ws.Cells.Item(26, 28).Value = 0.95
$ws.Cells.Item(26, 29).Value = -1
$ws.Cells.Item(27, 2).Value = 5143712
$ws.Cells.Item(27, 6).Value = 'Lyngby'
$ws.Cells.Item(27, 7).Value = 'AC Horsens'
$ws.Cells.Item(27, 8).Value = 1
$ws.Cells.Item(27, 9).Value = 1
$ws.Cells.Item(27, 11).Value = 2.2
$ws.Cells.Item(27, 12).Value = 3.4
$ws.Cells.Item(27, 13).Value = 3.2
$ws.Cells.Item(27, 14).Value = 2.1
$ws.Cells.Item(27, 15).Value = 3.4
$ws.Cells.Item(27, 16).Value = 3.75
$ws.Cells.Item(27, 18).Value = 1.8
$ws.Cells.Item(27, 19).Value = 2.05
$ws.Cells.Item(27, 20).Value = 2.25
$ws.Cells.Item(27, 21).Value = 1.85
$ws.Cells.Item(27, 22).Value = 2
$ws.Cells.Item(27, 24).Value = 2.4
$ws.Cells.Item(27, 27).Value = 0.5249999999999999
$ws.Cells.Item(27, 28).Value = -0.5
$ws.Cells.Item(27, 29).Value = 0.5
$ws.Cells.Item(29, 7).Value = 'Odense BK'
$ws.Cells.Item(32, 6).Value = 'Odense BK'
$ws.Cells.Item(33, 6).Value = 'Lyngby'
$ws.Cells.Item(39, 7).Value = 'Odense BK'
$ws.Cells.Item(40, 7).Value = 'Lyngby'
$ws.Cells.Item(45, 6).Value = 'Odense BK'
$ws.Cells.Item(46, 6).Value = 'Lyngby'
$ws.Cells.Item(55, 6).Value = 'Odense BK'
$ws.Cells.Item(55, 7).Value = 'Lyngby'
$ws.Cells.Item(56, 7).Value = 'Odense BK'
$ws.Cells.Item(57, 7).Value = 'Lyngby'
$ws.Cells.Item(63, 6).Value = 'Lyngby'
$ws.Cells.Item(64, 7).Value = 'Odense BK'
$ws.Cells.Item(68, 6).Value = 'Odense BK'
$ws.Cells.Item(69, 7).Value = 'Lyngby'
$ws.Cells.Item(75, 6).Value = 'Lyngby'
$ws.Cells.Item(75, 7).Value = 'Odense BK'
$ws.Cells.Item(80, 6).Value = 'Odense BK'
$ws.Cells.Item(81, 2).Value = 6445255
$ws.Cells.Item(81, 6).Value = 'Silkeborg IF'
$ws.Cells.Item(81, 7).Value = 'Midtjylland'
$ws.Cells.Item(81, 8).Value = 3
$ws.Cells.Item(81, 9).Value = 3
$ws.Cells.Item(81, 10).Value = 'D'
$ws.Cells.Item(81, 11).Value = 2.8
$ws.Cells.Item(81, 12).Value = 3.75
$ws.Cells.Item(81, 13).Value = 2.2
$ws.Cells.Item(81, 14).Value = 4
$ws.Cells.Item(81, 16).Value = 1.85
$ws.Cells.Item(81, 17).Value = 0.5
$ws.Cells.Item(81, 18).Value = 2.05
$ws.Cells.Item(81, 19).Value = 1.85
$ws.Cells.Item(81, 20).Value = 2.75
$ws.Cells.Item(81, 21).Value = 1.85
$ws.Cells.Item(81, 22).Value = 2
$ws.Cells.Item(81, 23).Value = -1
$ws.Cells.Item(81, 24).Value = 2.75
$ws.Cells.Item(81, 26).Value = 1.05
$ws.Cells.Item(81, 28).Value = 0.8500000000000001
$ws.Cells.Item(82, 2).Value = 6478389
$ws.Cells.Item(82, 6).Value = 'Lyngby'
$ws.Cells.Item(82, 7).Value = 'AaB'
$ws.Cells.Item(82, 8).Value = 2
$ws.Cells.Item(82, 9).Value = 1
$ws.Cells.Item(82, 10).Value = 'H'
$ws.Cells.Item(82, 11).Value = 2.75
$ws.Cells.Item(82, 12).Value = 3.7
$ws.Cells.Item(82, 13).Value = 2.3
$ws.Cells.Item(82, 14).Value = 3.2
$ws.Cells.Item(82, 16).Value = 2.1
$ws.Cells.Item(82, 17).Value = 0.25
$ws.Cells.Item(82, 18).Value = 2.025
$ws.Cells.Item(82, 19).Value = 1.825
$ws.Cells.Item(82, 20).Value = 2.5
$ws.Cells.Item(82, 21).Value = 1.875
$ws.Cells.Item(82, 22).Value = 1.975
$ws.Cells.Item(82, 23).Value = 2.2
$ws.Cells.Item(82, 24).Value = -1
$ws.Cells.Item(82, 26).Value = 1.025
$ws.Cells.Item(82, 28).Value = 0.875
$ws.Cells.Item(86, 7).Value = 'Odense BK'
$ws.Cells.Item(87, 7).Value = 'Lyngby'
$ws.Cells.Item(94, 6).Value = 'Lyngby'
$ws.Cells.Item(96, 6).Value = 'Odense BK'
$ws.Cells.Item(99, 7).Value = 'Lyngby'
$ws.Cells.Item(103, 7).Value = 'Odense BK'
$ws.Cells.Item(108, 6).Value = 'Lyngby'
$ws.Cells.Item(110, 6).Value = 'Odense BK'
$ws.Cells.Item(111, 7).Value = 'Odense BK'
$ws.Cells.Item(114, 7).Value = 'Lyngby'
$ws.Cells.Item(119, 6).Value = 'Lyngby'
$ws.Cells.Item(121, 6).Value = 'Odense BK'
$ws.Cells.Item(125, 7).Value = 'Lyngby'
$ws.Cells.Item(126, 7).Value = 'Odense BK'
$ws.Cells.Item(129, 6).Value = 'Odense BK'
$ws.Cells.Item(131, 6).Value = 'Lyngby'
$ws.Cells.Item(138, 7).Value = 'Lyngby'
$ws.Cells.Item(140, 6).Value = 'Odense BK'
$ws.Cells.Item(141, 6).Value = 'Lyngby'
$ws.Cells.Item(145, 7).Value = 'Odense BK'
$ws.Cells.Item(152, 6).Value = 'Odense BK'
$ws.Cells.Item(152, 7).Value = 'Lyngby'
$ws.Cells.Item(153, 7).Value = 'Lyngby'
$ws.Cells.Item(155, 7).Value = 'Odense BK'
$ws.Cells.Item(161, 6).Value = 'Lyngby'
$ws.Cells.Item(164, 6).Value = 'Odense BK'
$ws.Cells.Item(165, 7).Value = 'Lyngby'
$ws.Cells.Item(168, 7).Value = 'Odense BK'
$ws.Cells.Item(171, 6).Value = 'Lyngby'
$ws.Cells.Item(171, 7).Value = 'Odense BK'
$ws.Cells.Item(179, 7).Value = 'Lyngby'
$ws.Cells.Item(181, 6).Value = 'Odense BK'
$ws.Cells.Item(186, 7).Value = 'Odense BK'
$ws.Cells.Item(187, 6).Value = 'Lyngby'
$ws.Cells.Item(190, 2).Value = 6779673
$ws.Cells.Item(190, 7).Value = 'Silkeborg IF'
$ws.Cells.Item(190, 8).Value = 2
$ws.Cells.Item(190, 9).Value = 0
$ws.Cells.Item(190, 10).Value = 'H'
$ws.Cells.Item(190, 11).Value = 2.9
$ws.Cells.Item(190, 12).Value = 3.5
$ws.Cells.Item(190, 13).Value = 2.2
$ws.Cells.Item(190, 14).Value = 3.1
$ws.Cells.Item(190, 15).Value = 3.4
$ws.Cells.Item(190, 16).Value = 2.3
$ws.Cells.Item(190, 17).Value = 0.25
$ws.Cells.Item(190, 18).Value = 1.9
$ws.Cells.Item(190, 19).Value = 2
$ws.Cells.Item(190, 21).Value = 1.925
$ws.Cells.Item(190, 22).Value = 1.925
$ws.Cells.Item(190, 23).Value = 2.1
$ws.Cells.Item(190, 24).Value = -1
$ws.Cells.Item(190, 26).Value = 0.8999999999999999
$ws.Cells.Item(190, 29).Value = 0.925
$ws.Cells.Item(191, 2).Value = 6779676
$ws.Cells.Item(191, 7).Value = 'FC Nordsjaelland'
$ws.Cells.Item(191, 8).Value = 1
$ws.Cells.Item(191, 9).Value = 1
$ws.Cells.Item(191, 10).Value = 'D'
$ws.Cells.Item(191, 11).Value = 4.333
$ws.Cells.Item(191, 12).Value = 3.6
$ws.Cells.Item(191, 13).Value = 1.75
$ws.Cells.Item(191, 14).Value = 4.2
$ws.Cells.Item(191, 15).Value = 3.5
$ws.Cells.Item(191, 16).Value = 1.909
$ws.Cells.Item(191, 17).Value = 0.5
$ws.Cells.Item(191, 18).Value = 1.925
$ws.Cells.Item(191, 19).Value = 1.925
$ws.Cells.Item(191, 21).Value = 1.9
$ws.Cells.Item(191, 22).Value = 1.95
$ws.Cells.Item(191, 23).Value = -1
$ws.Cells.Item(191, 24).Value = 2.5
$ws.Cells.Item(191, 26).Value = 0.925
$ws.Cells.Item(191, 29).Value = 0.95
